$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds weekly price rows. Two new rows of data (Peru origin melons)
# are inserted before the existing row 53, pushing the former rows 53-64 down
# to rows 55-66 (dimension grows from A1:R64 to A1:R66).
$ws.Rows("53:54").Insert()

# --- New row 53 (Tuna / Primera, Peru, week of 2021-04-26) ---
$ws.Range("A53").Value = 12
$ws.Range("B53").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44312
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 100112027
$ws.Range("G53").Value = "Melón"
$ws.Range("H53").Value = "Tuna"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 1500
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = 1500
$ws.Range("N53").Value = "$/unidad"
$ws.Range("O53").Value = "Perú"
$ws.Range("P53").Value = 1500
$ws.Range("Q53").Value = 1
$ws.Range("R53").Value = "Hortaliza"

# --- New row 54 (Tuna / Primera, Peru, week of 2021-04-19) ---
$ws.Range("A54").Value = 12
$ws.Range("B54").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C54").Value = "Metropolitana"
$ws.Range("D54").Value = 44305
$ws.Range("E54").Value = 13
$ws.Range("F54").Value = 100112027
$ws.Range("G54").Value = "Melón"
$ws.Range("H54").Value = "Tuna"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 120
$ws.Range("K54").Value = 1500
$ws.Range("L54").Value = 1500
$ws.Range("M54").Value = 1500
$ws.Range("N54").Value = "$/unidad"
$ws.Range("O54").Value = "Perú"
$ws.Range("P54").Value = 1500
$ws.Range("Q54").Value = 1
$ws.Range("R54").Value = "Hortaliza"

# Column D holds dates; make sure the new rows carry the same date number
# format as the rest of the column (style index 2, already inherited from
# the insert, but set explicitly for safety).
$ws.Range("D53:D54").NumberFormat = $ws.Range("D55").NumberFormat
